# Weekly update: insert a new price record (row 34) for
# "Vega Monumental Concepción" / Ají / Inferno, shifting the existing
# rows 34-43 down to rows 35-44.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 34; Excel shifts rows 34:43 down to 35:44
# and extends the used range / dimension to A1:R44 automatically.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new record.
$ws.Range("A34").Value = 11
$ws.Range("B34").Value = "Vega Monumental Concepción"
$ws.Range("C34").Value = "Bíobío"
$ws.Range("D34").Value = 44463
$ws.Range("E34").Value = 8
$ws.Range("F34").Value = 100112021
$ws.Range("G34").Value = "Ají"
$ws.Range("H34").Value = "Inferno"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 30
$ws.Range("K34").Value = 41000
$ws.Range("L34").Value = 42000
$ws.Range("M34").Value = 41667
$ws.Range("N34").Value = "`$/caja 12 kilos"
$ws.Range("O34").Value = "Región de Arica y Parinacota"
$ws.Range("P34").Value = 3472
$ws.Range("Q34").Value = 12
$ws.Range("R34").Value = "Hortaliza"
